$d = $word.ActiveDocument
$d.Content.Find.Execute("At first this did reach the target accuracy", $true, $false, $false, $false, $false,
                         $true, 1, $false, "At first this did not reach the target accuracy", 2)
